$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: merge the three runs (with proofErr gramStart/gramEnd markers
# around "GPU") that make up the "Analyzing the Energy-Efficiency..." title
# into a single plain run with the full text.
# ---------------------------------------------------------------------------
$old1 = "Analyzing the Energy-Efficiency of Vision Kernels on Embedded CPU, GPU and FPGA Platforms"
$found1 = $d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $old1, 2)
if (-not $found1) {
    Write-Host "WARNING: change 1 text not found"
}

# ---------------------------------------------------------------------------
# Change 2: merge the three runs "{b1" / "8" / "}" into a single run "{b18}".
# ---------------------------------------------------------------------------
$found2 = $d.Content.Find.Execute("{b18}", $false, $false, $false, $false, $false, $true, 1, $false, "{b18}", 2)
if (-not $found2) {
    Write-Host "WARNING: change 2 text not found"
}

# ---------------------------------------------------------------------------
# Change 3: replace the empty ListParagraph-styled paragraph right after the
# "High-Throughput Line Buffer..." reference block with:
#   - a bare empty paragraph
#   - a bold heading-style paragraph "A Survey of Convolutional Neural
#     Networks on Edge with Reconfigurable Computing"
#   - a "{b21}" list paragraph
#   - a list paragraph whose whole text is a hyperlink to the google
#     redirect URL
# ---------------------------------------------------------------------------

# Locate the target empty paragraph: it is the ListParagraph-styled, empty
# paragraph that immediately follows the paragraph whose text is the mdpi
# hyperlink "https://www.mdpi.com/2313-433X/5/3/34".
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text
    if ($txt -like "*https://www.mdpi.com/2313-433X/5/3/34*") {
        $targetIndex = $i + 1
        break
    }
}

if ($targetIndex -eq -1) {
    Write-Host "ERROR: could not locate anchor paragraph for change 3"
} else {
    $targetPara = $d.Paragraphs($targetIndex)
    $targetRange = $d.Range($targetPara.Range.Start, $targetPara.Range.End)

    $placeholder = "@@HYPERLINK_PLACEHOLDER@@"

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body>' +
        '<w:p/>' +
        '<w:p><w:pPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:b/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:b/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>A Survey of Convolutional Neural Networks on Edge</w:t></w:r>' +
        '<w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:b/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:b/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>with Reconfigurable Computing</w:t></w:r>' +
        '</w:p>' +
        '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr></w:pPr>' +
        '<w:r><w:t>{b</w:t></w:r><w:r><w:t>21</w:t></w:r><w:r><w:t>}</w:t></w:r>' +
        '</w:p>' +
        '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr></w:pPr>' +
        '<w:r><w:t>' + $placeholder + '</w:t></w:r>' +
        '</w:p>' +
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $targetRange.InsertXML($xml)

    # Find the placeholder paragraph and turn its text into a real hyperlink.
    $linkParaIndex = -1
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $txt = $d.Paragraphs($i).Range.Text
        if ($txt -like "*$placeholder*") {
            $linkParaIndex = $i
            break
        }
    }

    if ($linkParaIndex -eq -1) {
        Write-Host "ERROR: could not find hyperlink placeholder paragraph"
    } else {
        $linkPara = $d.Paragraphs($linkParaIndex)
        $linkRange = $d.Range($linkPara.Range.Start, $linkPara.Range.End - 1)
        $url = "https://www.google.com/url?sa=t&rct=j&q=&esrc=s&source=web&cd=&ved=2ahUKEwiZpqvvyuTuAhWbbc0KHfaxCJwQFjAAegQIAxAC&url=https%3A%2F%2Fwww.mdpi.com%2F1999-4893%2F12%2F8%2F154%2Fpdf&usg=AOvVaw3SkNF_WK6DsxB2z7HCE3t-"
        $null = $d.Hyperlinks.Add($linkRange, $url, "", "", $url)
    }
}
